$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.945.32"
$ws.Range("E2").Value = "  +2.12%  "

$ws.Range("D3").Value = "3.585.39"
$ws.Range("E3").Value = "  -0.70%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.64"
$ws.Range("E5").Value = "  +3.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "653.55"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("E7").Value = "  +14.75%  "

$ws.Range("E8").Value = "  +3.80%  "

$ws.Range("E9").Value = "  +7.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  -0.06%  "

$ws.Range("D11").Value = "3.585.18"
$ws.Range("E11").Value = "  -0.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.57"
$ws.Range("E12").Value = "  +4.08%  "

$ws.Range("E13").Value = "  +1.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.43"
$ws.Range("E14").Value = "  +0.26%  "

$ws.Range("D15").Value = "4.249.00"
$ws.Range("E15").Value = "  -1.26%  "

$ws.Range("D16").Value = "96.667.32"
$ws.Range("E16").Value = "  +1.92%  "

$ws.Range("E17").Value = "  +2.53%  "

$ws.Range("D18").Value = "3.572.53"
$ws.Range("E18").Value = "  -1.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.77"
$ws.Range("E19").Value = "  -1.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.65"
$ws.Range("E20").Value = "  -1.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.04"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.535"
$ws.Range("E22").Value = "  +12.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "509.11"
$ws.Range("E23").Value = "  +1.70%  "

$ws.Range("E24").Value = "  -2.84%  "

$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000201"
$ws.Range("E25").Value = "  +3.84%  "

$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.94"
$ws.Range("E26").Value = "  +5.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "97.06"
$ws.Range("E27").Value = "  +2.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.15"
$ws.Range("E28").Value = "  +6.13%  "

$ws.Range("D29").Value = "3.777.05"
$ws.Range("E29").Value = "  -0.78%  "

$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.153"
$ws.Range("E30").Value = "  +11.60%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.04"
$ws.Range("E31").Value = "  -1.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.52"
$ws.Range("E32").Value = "  +3.29%  "

$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("E34").Value = "  +5.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  -0.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.44"
$ws.Range("E36").Value = "  -2.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.95"
$ws.Range("E37").Value = "  +12.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "626.07"
$ws.Range("E38").Value = "  +11.53%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.571"
$ws.Range("E39").Value = "  +2.94%  "

$ws.Range("E40").Value = "  +12.61%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.151"
$ws.Range("E41").Value = "  +1.54%  "

$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.914"
$ws.Range("E43").Value = "  +0.36%  "

$ws.Range("E44").Value = "  +6.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.84"
$ws.Range("E45").Value = "  +4.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0434"
$ws.Range("E46").Value = "  +5.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.31"
$ws.Range("E47").Value = "  +3.77%  "

$ws.Range("E48").Value = "  -0.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.06"
$ws.Range("E49").Value = "  -9.73%  "

$ws.Range("E50").Value = "  +0.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.32"
$ws.Range("E51").Value = "  +4.62%  "

